$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells we touch keep their text (string) representation
# instead of being auto-converted to numbers by Excel when values look numeric.
$priceCells = @("D2","D3","D5","D6","D7","D9","D11","D13","D14","D15","D16","D17","D21","D22","D23","D24","D25","D27","D29","D31","D32","D33","D34","D35","D37","D41","D45","D46","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update Coin / Link for the two rows that swapped ranking position (50 and 51)
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'

# Update Price column (D) values
$ws.Range('D2').Value = '61.473.39'
$ws.Range('D3').Value = '3.376.27'
$ws.Range('D5').Value = '406.95'
$ws.Range('D6').Value = '135.13'
$ws.Range('D7').Value = '0.593'
$ws.Range('D9').Value = '0.670'
$ws.Range('D11').Value = '42.60'
$ws.Range('D13').Value = '3.898.16'
$ws.Range('D14').Value = '8.42'
$ws.Range('D15').Value = '19.70'
$ws.Range('D16').Value = '3.360.66'
$ws.Range('D17').Value = '61.447.53'
$ws.Range('D21').Value = '3.20'
$ws.Range('D22').Value = '85.08'
$ws.Range('D23').Value = '314.18'
$ws.Range('D24').Value = '12.82'
$ws.Range('D25').Value = '3.14'
$ws.Range('D27').Value = '8.37'
$ws.Range('D29').Value = '7.63'
$ws.Range('D31').Value = '0.171'
$ws.Range('D32').Value = '2.58'
$ws.Range('D33').Value = '11.35'
$ws.Range('D34').Value = '0.999'
$ws.Range('D35').Value = '40.68'
$ws.Range('D37').Value = '51.91'
$ws.Range('D41').Value = '138.68'
$ws.Range('D45').Value = '4.04'
$ws.Range('D46').Value = '16.75'
$ws.Range('D48').Value = '21.29'
$ws.Range('D49').Value = '2.125.02'
$ws.Range('D50').Value = '2.29'
$ws.Range('D51').Value = '1.93'

# Update Volume(1h) column (E) values
$ws.Range('E2').Value = '  -0.98%  '
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E5').Value = '  -1.85%  '
$ws.Range('E6').Value = '  +7.79%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  +1.63%  '
$ws.Range('E10').Value = '  -5.05%  '
$ws.Range('E11').Value = '  +2.54%  '
$ws.Range('E12').Value = '  -1.04%  '
$ws.Range('E13').Value = '  -2.56%  '
$ws.Range('E14').Value = '  -1.19%  '
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('E16').Value = '  -2.59%  '
$ws.Range('E17').Value = '  -0.98%  '
$ws.Range('E18').Value = '  -1.72%  '
$ws.Range('E19').Value = '  -0.52%  '
$ws.Range('E20').Value = '  -3.50%  '
$ws.Range('E21').Value = '  -4.01%  '
$ws.Range('E22').Value = '  +3.54%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('E24').Value = '  -1.36%  '
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('E26').Value = '  +11.60%  '
$ws.Range('E27').Value = '  +6.12%  '
$ws.Range('E28').Value = '  -4.76%  '
$ws.Range('E29').Value = '  -2.82%  '
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('E31').Value = '  -1.47%  '
$ws.Range('E32').Value = '  +1.00%  '
$ws.Range('E33').Value = '  -1.91%  '
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('E35').Value = '  -3.30%  '
$ws.Range('E36').Value = '  -0.58%  '
$ws.Range('E37').Value = '  -0.50%  '
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('E39').Value = '  -2.62%  '
$ws.Range('E40').Value = '  -2.22%  '
$ws.Range('E41').Value = '  +3.02%  '
$ws.Range('E42').Value = '  -1.97%  '
$ws.Range('E43').Value = '  -1.04%  '
$ws.Range('E44').Value = '  +3.36%  '
$ws.Range('E45').Value = '  +3.52%  '
$ws.Range('E46').Value = '  -2.88%  '
$ws.Range('E47').Value = '  +1.21%  '
$ws.Range('E48').Value = '  -4.52%  '
$ws.Range('E49').Value = '  -3.86%  '
$ws.Range('E50').Value = '  -5.11%  '
$ws.Range('E51').Value = '  +1.75%  '
